# "added ltc4359 ckt + components"
#
# The original sheet had columns: A Board | B Component | C Part | D Link | E Quantity | F Comments
# The edit inserts a new "Subsystem" column after "Board", and appends six new
# component rows (LTC4359 Oring controller circuit + MOSFETs/TVS/Zener diodes)
# into the previously-blank filler rows 11-16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column B ("Subsystem") - shifts old B..F to C..G.
# ---------------------------------------------------------------------------
$ws.Columns("B").Insert()

# ---------------------------------------------------------------------------
# 2. Header row.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Subsystem"

# ---------------------------------------------------------------------------
# 3. Fill in the new Subsystem column for the existing power-supply rows.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "3.3V buck"
$ws.Range("B3").Value = "3.3V buck/5V buck"
$ws.Range("B4").Value = "3.3V buck/5V buck"
$ws.Range("B5").Value = "3.3V buck/5V buck"
$ws.Range("B6").Value = "3.3V buck/5V buck"
$ws.Range("B7").Value = "5V buck"
# Rows 8-10 (resistor / AC-DC converter / power monitor) intentionally keep no
# Subsystem value, matching the source data.

# ---------------------------------------------------------------------------
# 4. New components for the "Switching network" subsystem (Power Oring
#    controller circuit, MOSFETs, TVS diodes, Zener diode) in rows 11-16,
#    which were previously blank filler rows.
# ---------------------------------------------------------------------------
$newRows = @(
  @{ Row = 11; Subsystem = "Switching network"; Component = "Power Oring controller"; Part = "LTC4359CMS8#TRPBF"; Link = "https://www.digikey.com/en/products/detail/analog-devices-inc/LTC4359CMS8-TRPBF/3306822" },
  @{ Row = 12; Subsystem = "Switching network"; Component = "MOSFET";                 Part = "NTBGS1D5N06C";       Link = "https://www.digikey.com/en/products/detail/onsemi/NTBGS1D5N06C/14005228" },
  @{ Row = 13; Subsystem = "Switching network"; Component = "TVS diode";              Part = "SMAJ24CA-13-F";      Link = "https://www.digikey.com/en/products/detail/diodes-incorporated/SMAJ24CA-13-F/775748" },
  @{ Row = 14; Subsystem = "Switching network"; Component = "TVS diode";              Part = "SMAJ58A";            Link = "https://www.digikey.com/en/products/detail/littelfuse-inc/SMAJ58A/762310" },
  @{ Row = 15; Subsystem = "Switching network"; Component = "MOSFET";                 Part = "FDB13AN06A0";        Link = "https://www.digikey.com/en/products/detail/onsemi/FDB13AN06A0/978478" },
  @{ Row = 16; Subsystem = "Switching network"; Component = "Zener diode";            Part = "BZX84C12VLYT116";    Link = "https://www.digikey.com/en/products/detail/rohm-semiconductor/BZX84C12VLYT116/14682680" }
)

foreach ($r in $newRows) {
  $row = $r.Row
  $ws.Range("A$row").Value = "Power"
  $ws.Range("B$row").Value = $r.Subsystem
  $ws.Range("C$row").Value = $r.Component
  $ws.Range("D$row").Value = $r.Part
  $ws.Range("E$row").Value = $r.Link
  $ws.Hyperlinks.Add($ws.Range("E$row"), $r.Link)
}

# Re-apply the existing "Link" cell format (Hyperlink style, already used on
# E2/E3/... after the column insert) onto the new hyperlink cells so they
# match the rest of the sheet instead of picking up a brand-new style.
$ws.Range("E2").Copy()
$ws.Range("E11:E16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# D16 keeps the quirky pre-existing "s=5" look from the source workbook
# (same format as the blank quantity cells in rows 3-6, i.e. column F there).
$ws.Range("F3").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = "BZX84C12VLYT116"
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5. Column widths (best-effort autofit to match the resized columns).
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 19.45
$ws.Columns("C").ColumnWidth = 26.45
$ws.Columns("D").ColumnWidth = 29.88
$ws.Columns("E").ColumnWidth = 111.59
$ws.Columns("F").ColumnWidth = 10.74
$ws.Columns("G").ColumnWidth = 35.59

# ---------------------------------------------------------------------------
# 6. Selection / view state.
# ---------------------------------------------------------------------------
$ws.Range("D22").Select()
